$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DATASET column (A) values for rows 2-65 (row 6 unchanged)
$ws.Cells.Item(2, 1).Value = 9.71
$ws.Cells.Item(3, 1).Value = 9.83
$ws.Cells.Item(4, 1).Value = 10.08
$ws.Cells.Item(5, 1).Value = 11.32
$ws.Cells.Item(7, 1).Value = 13.52
$ws.Cells.Item(8, 1).Value = 14.12
$ws.Cells.Item(9, 1).Value = 15.89
$ws.Cells.Item(10, 1).Value = 16.09
$ws.Cells.Item(11, 1).Value = 17.34
$ws.Cells.Item(12, 1).Value = 18.48
$ws.Cells.Item(13, 1).Value = 18.52
$ws.Cells.Item(14, 1).Value = 18.67
$ws.Cells.Item(15, 1).Value = 19.05
$ws.Cells.Item(16, 1).Value = 19.13
$ws.Cells.Item(17, 1).Value = 20.14
$ws.Cells.Item(18, 1).Value = 21
$ws.Cells.Item(19, 1).Value = 21.17
$ws.Cells.Item(20, 1).Value = 21.83
$ws.Cells.Item(21, 1).Value = 22
$ws.Cells.Item(22, 1).Value = 22.61
$ws.Cells.Item(23, 1).Value = 22.89
$ws.Cells.Item(24, 1).Value = 23.08
$ws.Cells.Item(25, 1).Value = 23.38
$ws.Cells.Item(26, 1).Value = 23.42
$ws.Cells.Item(27, 1).Value = 23.7
$ws.Cells.Item(28, 1).Value = 23.8
$ws.Cells.Item(29, 1).Value = 24.31
$ws.Cells.Item(30, 1).Value = 24.62
$ws.Cells.Item(31, 1).Value = 24.73
$ws.Cells.Item(32, 1).Value = 24.88
$ws.Cells.Item(33, 1).Value = 25.08
$ws.Cells.Item(34, 1).Value = 25.12
$ws.Cells.Item(35, 1).Value = 25.21
$ws.Cells.Item(36, 1).Value = 25.64
$ws.Cells.Item(37, 1).Value = 26.01
$ws.Cells.Item(38, 1).Value = 26.17
$ws.Cells.Item(39, 1).Value = 26.21
$ws.Cells.Item(40, 1).Value = 26.41
$ws.Cells.Item(41, 1).Value = 26.87
$ws.Cells.Item(42, 1).Value = 27.13
$ws.Cells.Item(43, 1).Value = 27.33
$ws.Cells.Item(44, 1).Value = 27.42
$ws.Cells.Item(45, 1).Value = 27.53
$ws.Cells.Item(46, 1).Value = 27.91
$ws.Cells.Item(47, 1).Value = 28.09
$ws.Cells.Item(48, 1).Value = 28.18
$ws.Cells.Item(49, 1).Value = 28.31
$ws.Cells.Item(50, 1).Value = 28.42
$ws.Cells.Item(51, 1).Value = 28.48
$ws.Cells.Item(52, 1).Value = 29
$ws.Cells.Item(53, 1).Value = 30.32
$ws.Cells.Item(54, 1).Value = 31.12
$ws.Cells.Item(55, 1).Value = 32.33
$ws.Cells.Item(56, 1).Value = 32.33
$ws.Cells.Item(57, 1).Value = 33.45
$ws.Cells.Item(58, 1).Value = 34.42
$ws.Cells.Item(59, 1).Value = 35.89
$ws.Cells.Item(60, 1).Value = 36.7
$ws.Cells.Item(61, 1).Value = 38.45
$ws.Cells.Item(62, 1).Value = 39.08
$ws.Cells.Item(63, 1).Value = 40.08
$ws.Cells.Item(64, 1).Value = 41.76
$ws.Cells.Item(65, 1).Value = 42.67

# Clear former dataset tail (rows 66-101) - dataset now ends at row 65
$ws.Range("A66:A101").ClearContents()

# Re-sort the dataset range to match the new extent
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2"))
$ws.Sort.SetRange($ws.Range("A2:A65"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Update active selection
$ws.Range("D7").Select() | Out-Null